$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.044933333333333
$ws.Range("H2").Value = 18.1348
$ws.Range("I2").Value = 0.9708761253868625
$ws.Range("J2").Value = 0.9708761253868624
$ws.Range("M2").Value = 11.651608
$ws.Range("N2").Value = 34.954824
$ws.Range("O2").Value = 0.1892813629236475
$ws.Range("P2").Value = 0.1892813629236474
$ws.Range("Q2").Value = 70.43319358613334
$ws.Range("R2").Value = 633.8987422752
$ws.Range("S2").Value = 0.1837687562432554
$ws.Range("T2").Value = 0.1837687562432553
# Row 3
$ws.Range("G3").Value = 6.044933333333333
$ws.Range("H3").Value = 18.1348
$ws.Range("I3").Value = 0.9708761253868625
$ws.Range("J3").Value = 0.9708761253868624
$ws.Range("O3").Value = 0.4419371310876561
$ws.Range("P3").Value = 0.4419371310876561
$ws.Range("Q3").Value = 164.4485385460444
$ws.Range("R3").Value = 1480.0368469144
$ws.Range("S3").Value = 0.4290662094949695
$ws.Range("T3").Value = 0.4290662094949694
# Row 4
$ws.Range("G4").Value = 6.044933333333333
$ws.Range("H4").Value = 18.1348
$ws.Range("I4").Value = 0.9708761253868625
$ws.Range("J4").Value = 0.9708761253868624
$ws.Range("M4").Value = 8.657178999999999
$ws.Range("N4").Value = 25.971537
$ws.Range("O4").Value = 0.1406366091439035
$ws.Range("P4").Value = 0.1406366091439035
$ws.Range("Q4").Value = 52.33206990973333
$ws.Range("R4").Value = 470.9886291875999
$ws.Range("S4").Value = 0.1365407261731796
$ws.Range("T4").Value = 0.1365407261731796
# Row 5
$ws.Range("G5").Value = 6.044933333333333
$ws.Range("H5").Value = 18.1348
$ws.Range("I5").Value = 0.9708761253868625
$ws.Range("J5").Value = 0.9708761253868624
$ws.Range("M5").Value = 5.488499666666667
$ws.Range("N5").Value = 16.465499
$ws.Range("O5").Value = 0.08916114387925267
$ws.Range("P5").Value = 0.08916114387925267
$ws.Range("Q5").Value = 33.17761458502222
$ws.Range("R5").Value = 298.5985312652
$ws.Range("S5").Value = 0.08656442590454941
$ws.Range("T5").Value = 0.08656442590454939
# Row 6
$ws.Range("G6").Value = 6.044933333333333
$ws.Range("H6").Value = 18.1348
$ws.Range("I6").Value = 0.9708761253868625
$ws.Range("J6").Value = 0.9708761253868624
$ws.Range("M6").Value = 4.091608333333333
$ws.Range("N6").Value = 12.274825
$ws.Range("O6").Value = 0.06646852536431769
$ws.Range("P6").Value = 0.06646852536431769
$ws.Range("Q6").Value = 24.73349960111111
$ws.Range("R6").Value = 222.60149641
$ws.Range("S6").Value = 0.06453270436588715
$ws.Range("T6").Value = 0.06453270436588715
# Row 7
$ws.Range("G7").Value = 6.044933333333333
$ws.Range("H7").Value = 18.1348
$ws.Range("I7").Value = 0.9708761253868625
$ws.Range("J7").Value = 0.9708761253868624
$ws.Range("M7").Value = 4.463825666666667
$ws.Range("N7").Value = 13.391477
$ws.Range("O7").Value = 0.07251522760122259
$ws.Range("P7").Value = 0.07251522760122257
$ws.Range("Q7").Value = 26.98352856662222
$ws.Range("R7").Value = 242.8517570996
$ws.Range("S7").Value = 0.07040330320502146
$ws.Range("T7").Value = 0.07040330320502143
# Row 8
$ws.Range("I8").Value = 0.01821359071319307
$ws.Range("J8").Value = 0.01821359071319307
$ws.Range("M8").Value = 11.651608
$ws.Range("N8").Value = 34.954824
$ws.Range("O8").Value = 0.1892813629236475
$ws.Range("P8").Value = 0.1892813629236474
$ws.Range("Q8").Value = 1.321323418154667
$ws.Range("R8").Value = 11.891910763392
$ws.Range("S8").Value = 0.003447493273926673
$ws.Range("T8").Value = 0.003447493273926672
# Row 9
$ws.Range("I9").Value = 0.01821359071319307
$ws.Range("J9").Value = 0.01821359071319307
$ws.Range("O9").Value = 0.4419371310876561
$ws.Range("P9").Value = 0.4419371310876561
$ws.Range("S9").Value = 0.008049262026593323
$ws.Range("T9").Value = 0.008049262026593323
# Row 10
$ws.Range("I10").Value = 0.01821359071319307
$ws.Range("J10").Value = 0.01821359071319307
$ws.Range("M10").Value = 8.657178999999999
$ws.Range("N10").Value = 25.971537
$ws.Range("O10").Value = 0.1406366091439035
$ws.Range("P10").Value = 0.1406366091439035
$ws.Range("Q10").Value = 0.9817471844106666
$ws.Range("R10").Value = 8.835724659696
$ws.Range("S10").Value = 0.002561497638238365
$ws.Range("T10").Value = 0.002561497638238365
# Row 11
$ws.Range("I11").Value = 0.01821359071319307
$ws.Range("J11").Value = 0.01821359071319307
$ws.Range("M11").Value = 5.488499666666667
$ws.Range("N11").Value = 16.465499
$ws.Range("O11").Value = 0.08916114387925267
$ws.Range("P11").Value = 0.08916114387925267
$ws.Range("Q11").Value = 0.6224104981991111
$ws.Range("R11").Value = 5.601694483792
$ws.Range("S11").Value = 0.001623944582136828
$ws.Range("T11").Value = 0.001623944582136828
# Row 12
$ws.Range("I12").Value = 0.01821359071319307
$ws.Range("J12").Value = 0.01821359071319307
$ws.Range("M12").Value = 4.091608333333333
$ws.Range("N12").Value = 12.274825
$ws.Range("O12").Value = 0.06646852536431769
$ws.Range("P12").Value = 0.06646852536431769
$ws.Range("Q12").Value = 0.4639992959555556
$ws.Range("R12").Value = 4.1759936636
$ws.Range("S12").Value = 0.001210630516295175
$ws.Range("T12").Value = 0.001210630516295175
# Row 13
$ws.Range("I13").Value = 0.01821359071319307
$ws.Range("J13").Value = 0.01821359071319307
$ws.Range("M13").Value = 4.463825666666667
$ws.Range("N13").Value = 13.391477
$ws.Range("O13").Value = 0.07251522760122259
$ws.Range("P13").Value = 0.07251522760122257
$ws.Range("Q13").Value = 0.5062097341351111
$ws.Range("R13").Value = 4.555887607216
$ws.Range("S13").Value = 0.00132076267600271
$ws.Range("T13").Value = 0.001320762676002709
# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.06793033333333333
$ws.Range("H14").Value = 0.203791
$ws.Range("I14").Value = 0.01091028389994453
$ws.Range("J14").Value = 0.01091028389994453
$ws.Range("M14").Value = 11.651608
$ws.Range("N14").Value = 34.954824
$ws.Range("O14").Value = 0.1892813629236475
$ws.Range("P14").Value = 0.1892813629236474
$ws.Range("Q14").Value = 0.7914976153093334
$ws.Range("R14").Value = 7.123478537784001
$ws.Range("S14").Value = 0.002065113406465429
$ws.Range("T14").Value = 0.002065113406465429
# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.06793033333333333
$ws.Range("H15").Value = 0.203791
$ws.Range("I15").Value = 0.01091028389994453
$ws.Range("J15").Value = 0.01091028389994453
$ws.Range("O15").Value = 0.4419371310876561
$ws.Range("P15").Value = 0.4419371310876561
$ws.Range("Q15").Value = 1.848001197633111
$ws.Range("R15").Value = 16.632010778698
$ws.Range("S15").Value = 0.00482165956609333
$ws.Range("T15").Value = 0.00482165956609333
# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.06793033333333333
$ws.Range("H16").Value = 0.203791
$ws.Range("I16").Value = 0.01091028389994453
$ws.Range("J16").Value = 0.01091028389994453
$ws.Range("M16").Value = 8.657178999999999
$ws.Range("N16").Value = 25.971537
$ws.Range("O16").Value = 0.1406366091439035
$ws.Range("P16").Value = 0.1406366091439035
$ws.Range("Q16").Value = 0.5880850551963333
$ws.Range("R16").Value = 5.292765496766999
$ws.Range("S16").Value = 0.001534385332485522
$ws.Range("T16").Value = 0.001534385332485522
# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.06793033333333333
$ws.Range("H17").Value = 0.203791
$ws.Range("I17").Value = 0.01091028389994453
$ws.Range("J17").Value = 0.01091028389994453
$ws.Range("M17").Value = 5.488499666666667
$ws.Range("N17").Value = 16.465499
$ws.Range("O17").Value = 0.08916114387925267
$ws.Range("P17").Value = 0.08916114387925267
$ws.Range("Q17").Value = 0.3728356118565556
$ws.Range("R17").Value = 3.355520506709
$ws.Range("S17").Value = 0.0009727733925664483
$ws.Range("T17").Value = 0.0009727733925664483
# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.06793033333333333
$ws.Range("H18").Value = 0.203791
$ws.Range("I18").Value = 0.01091028389994453
$ws.Range("J18").Value = 0.01091028389994453
$ws.Range("M18").Value = 4.091608333333333
$ws.Range("N18").Value = 12.274825
$ws.Range("O18").Value = 0.06646852536431769
$ws.Range("P18").Value = 0.06646852536431769
$ws.Range("Q18").Value = 0.2779443179527777
$ws.Range("R18").Value = 2.501498861575
$ws.Range("S18").Value = 0.0007251904821353701
$ws.Range("T18").Value = 0.0007251904821353701
# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.06793033333333333
$ws.Range("H19").Value = 0.203791
$ws.Range("I19").Value = 0.01091028389994453
$ws.Range("J19").Value = 0.01091028389994453
$ws.Range("M19").Value = 4.463825666666667
$ws.Range("N19").Value = 13.391477
$ws.Range("O19").Value = 0.07251522760122259
$ws.Range("P19").Value = 0.07251522760122259
$ws.Range("Q19").Value = 0.3032291654785556
$ws.Range("R19").Value = 2.729062489307
$ws.Range("S19").Value = 0.0007911617201984321
$ws.Range("T19").Value = 0.000791161720198432

Write-Output "Done"